$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row (row 11): Right marks 5 -> 4, Wrong marks -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Total row (row 12): Right total 120 -> 96, and summary text updated accordingly
$ws.Range("B12").Value = 96
$ws.Range("E12").Value = "96 / 112"
